$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "65.736.49", "0.0611") that must stay
# literal text -- Excel auto-converts numeric-looking input to Number type,
# which would silently rewrite values / drop formatting (trailing zeros,
# multi-dot big numbers, subscript notation). Force Text format first so the
# assignments below round-trip byte-for-byte.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "65.736.49"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3
$ws.Range("D3").Value = "2.680.44"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "601.44"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").Value = "156.88"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").Value = "  +6.03%  "

# Row 9
$ws.Range("D9").Value = "0.130"
$ws.Range("E9").Value = "  +4.65%  "

# Row 10
$ws.Range("E10").Value = "  -0.76%  "

# Row 11
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  -3.41%  "

# Row 12
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "29.35"
$ws.Range("E13").Value = "  -2.90%  "

# Row 14
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("D15").Value = "3.161.15"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16
$ws.Range("D16").Value = "65.589.54"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "2.686.22"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
$ws.Range("D18").Value = "12.84"
$ws.Range("E18").Value = "  +1.26%  "

# Row 19
$ws.Range("E19").Value = "  -1.82%  "

# Row 20
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
$ws.Range("D21").Value = "352.07"
$ws.Range("E21").Value = "  -2.31%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "69.71"
$ws.Range("E23").Value = "  -0.73%  "

# Row 24
$ws.Range("E24").Value = "  +4.35%  "

# Row 25
$ws.Range("D25").Value = "9.67"
$ws.Range("E25").Value = "  -1.65%  "

# Row 27
$ws.Range("D27").Value = "0.168"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("E28").Value = "  -5.45%  "

# Row 29
$ws.Range("D29").Value = "8.09"
$ws.Range("E29").Value = "  -1.96%  "

# Row 30
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("E31").Value = "  -2.54%  "

# Row 32
$ws.Range("D32").Value = "531.41"
$ws.Range("E32").Value = "  +0.01%  "

# Row 33
$ws.Range("E33").Value = "  -2.35%  "

# Row 34
$ws.Range("D34").Value = "6.48"
$ws.Range("E34").Value = "  -2.38%  "

# Row 35
$ws.Range("E35").Value = "  +0.83%  "

# Row 36
$ws.Range("D36").Value = "0.425"
$ws.Range("E36").Value = "  -1.60%  "

# Row 37
$ws.Range("D37").Value = "20.53"
$ws.Range("E37").Value = "  -1.42%  "

# Row 38
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("D39").Value = "157.98"
$ws.Range("E39").Value = "  -3.09%  "

# Row 40
$ws.Range("E40").Value = "  -2.22%  "

# Row 41
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("D42").Value = "164.84"
$ws.Range("E42").Value = "  -2.46%  "

# Row 43
$ws.Range("D43").Value = "4.15"
$ws.Range("E43").Value = "  -0.60%  "

# Row 44
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +2.31%  "

# Row 45
$ws.Range("D45").Value = "0.0611"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").Value = "22.86"
$ws.Range("E46").Value = "  -2.13%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.643"
$ws.Range("E47").Value = "  -2.28%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0259"
$ws.Range("E48").Value = "  -2.53%  "

# Row 49
$ws.Range("D49").Value = "0.0₆0262"
$ws.Range("E49").Value = "  +15.07%  "

# Row 50
$ws.Range("E50").Value = "  +2.24%  "

# Row 51
$ws.Range("D51").Value = "20.13"
$ws.Range("E51").Value = "  -4.73%  "

